$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 16

$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value = "Bíobío"
$ws.Cells.Item($row, 4).Value = 44617
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100101
$ws.Cells.Item($row, 8).Value = "Berries"
$ws.Cells.Item($row, 9).Value = 100101004
$ws.Cells.Item($row, 10).Value = "Frambuesa"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 200
$ws.Cells.Item($row, 14).Value = 6000
$ws.Cells.Item($row, 15).Value = 7000
$ws.Cells.Item($row, 16).Value = 6500
$ws.Cells.Item($row, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item($row, 18).Value = "Región de Ñuble"
$ws.Cells.Item($row, 19).Value = 3250
$ws.Cells.Item($row, 20).Value = 2
